$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The regression gained two new explanatory series ("Crisis" / "Credit
# Allocation"): A_LF and FFR_LF. Insert two new columns between the
# existing A_C (col B) and LF_FFR (col C) columns, pushing LF_FFR to
# column E, then label the new header cells.
$ws.Range("C1:D1").EntireColumn.Insert()
$ws.Range("C1").Value = "A_LF"
$ws.Range("D1").Value = "FFR_LF"

# Refresh the regression output (params row 2, pvalue row 3) for all four
# explanatory columns now that the model has been re-estimated.
$ws.Range("B2").Value = 0.464722090331432
$ws.Range("C2").Value = -0.00961328564128364
$ws.Range("D2").Value = 2.378522371567055
$ws.Range("E2").Value = 0.4078882327616589

$ws.Range("B3").Value = 0.01503618942925478
$ws.Range("C3").Value = 0.003808133919186218
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

$wb.Save()
